$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking strings (e.g. "0.499")
# are not auto-converted to numbers by Excel, matching the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.658.44"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "1.635.42"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "213.10"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "0.499"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("D9").Value = "0.0623"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").Value = "19.14"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("E11").Value = "  +3.20%  "
$ws.Range("D12").Value = "1.864.13"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "1.642.30"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").Value = "26.678.31"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").Value = "63.26"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").Value = "217.35"
$ws.Range("E19").Value = "  +7.13%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "9.51"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").Value = "6.22"
$ws.Range("E23").Value = "  +2.62%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "149.21"
$ws.Range("E25").Value = "  +4.17%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "0.121"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").Value = "6.87"
$ws.Range("E28").Value = "  +4.49%  "
$ws.Range("D29").Value = "15.47"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").Value = "0.0518"
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  +3.25%  "
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("E35").Value = "  -2.34%  "
$ws.Range("D36").Value = "1.188.53"
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("E37").Value = "  +5.29%  "
$ws.Range("D38").Value = "0.808"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "0.506"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").Value = "5.41"
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").Value = "0.793"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "1.773.43"
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").Value = "92.36"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").Value = "54.75"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").Value = "7.65"
$ws.Range("E49").Value = "  +4.93%  "
$ws.Range("D50").Value = "0.411"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("E51").Value = "  +0.25%  "
